# Update column G ("K") values for rows 2-21 per the new save_data regen
# (regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 6
    3  = 2
    4  = 7
    5  = 3
    6  = 2
    7  = 5
    8  = 3
    9  = 3
    10 = 2
    11 = 3
    12 = 1
    13 = 5
    14 = 1
    15 = 1
    16 = 5
    17 = 1
    18 = 1
    19 = 3
    20 = 1
    21 = 4
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
